# Change log: #70 aggregation script / #69 hardware monitoring / #62 report draft
# Adds rows for exp_71 .. exp_107 (experiment overview sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fill in the learning-algorithm / estimator / game columns for the
#     trailing rows that previously only had an experiment name -----------
$ws.Range("B72").Value = "DQN"
$ws.Range("C72").Value = "Boosting"
$ws.Range("D72").Value = "CartPole"

$ws.Range("B73").Value = "DQN"
$ws.Range("C73").Value = "MC-Dropout"
$ws.Range("D73").Value = "CartPole"

$ws.Range("B74").Value = "DQN"
$ws.Range("C74").Value = "Boosting"
$ws.Range("D74").Value = "CartPole"
$ws.Range("E74").Value = "20 identical experiments`nused for evaluation for report`nit uses the same dropout as the MC evaluation and should therefore be comparable"
$ws.Range("E74").WrapText = $true
$ws.Rows.Item(74).RowHeight = 75

$ws.Range("B75").Value = "DQN"
$ws.Range("C75").Value = "Ensemble"
$ws.Range("D75").Value = "CartPole"
$ws.Range("E75").Value = "Series of 10 experiments"

$ws.Range("B76").Value = "DQN"
$ws.Range("C76").Value = "Ensemble"
$ws.Range("D76").Value = "CartPole"

# --- new experiment rows exp_76 .. exp_107 (rows 77-108), name only ------
$ws.Range("A77").Value = "exp_76"
$ws.Range("A78").Value = "exp_77"
$ws.Range("A79").Value = "exp_78"
$ws.Range("A80").Value = "exp_79"
$ws.Range("A81").Value = "exp_80"
$ws.Range("A82").Value = "exp_81"
$ws.Range("A83").Value = "exp_82"
$ws.Range("A84").Value = "exp_83"
$ws.Range("A85").Value = "exp_84"
$ws.Range("A86").Value = "exp_85"
$ws.Range("A87").Value = "exp_86"
$ws.Range("A88").Value = "exp_87"
$ws.Range("A89").Value = "exp_88"
$ws.Range("A90").Value = "exp_89"
$ws.Range("A91").Value = "exp_90"
$ws.Range("A92").Value = "exp_91"
$ws.Range("A93").Value = "exp_92"
$ws.Range("A94").Value = "exp_93"
$ws.Range("A95").Value = "exp_94"
$ws.Range("A96").Value = "exp_95"
$ws.Range("A97").Value = "exp_96"
$ws.Range("A98").Value = "exp_97"
$ws.Range("A99").Value = "exp_98"
$ws.Range("A100").Value = "exp_99"
$ws.Range("A101").Value = "exp_100"
$ws.Range("A102").Value = "exp_101"
$ws.Range("A103").Value = "exp_102"
$ws.Range("A104").Value = "exp_103"
$ws.Range("A105").Value = "exp_104"
$ws.Range("A106").Value = "exp_105"
$ws.Range("A107").Value = "exp_106"
$ws.Range("A108").Value = "exp_107"

# --- cosmetic: widen comment column, match the author's final selection --
$ws.Columns.Item(5).ColumnWidth = 38.5

[void]$ws.Range("B108").Select()
